$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighboring header cell (G1) onto the
# new header cell (H1) so it reuses the same cell style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column, plain numeric value like its neighbors.
$ws.Range("H2").Value = 0
